# Update the product name on the ProductLoanInput sheet: insert a dash
# after "199" so it reads "199-MS-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"
$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$wsInput.Range("B1").Value = "199-MS-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"
$wsOutput.Range("B1").Value = "199-MS-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

# Update the selection on the input sheet to B1
$wsInput.Range("B1").Select()

# Activate the output sheet and select B1 there, leaving it as the active tab
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
